# Added null check to GetAllProjectsForVertical
# Update the sample data row on Sheet1 to reflect a new project entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Saturday Works"
$ws.Range("D2").Value = "Bits 10"
$ws.Range("E2").Value = "New Landing page"
$ws.Range("F2").Value = "N/A"

# Move the active selection to G2, matching the author's final cursor position.
[void]$ws.Range("G2").Select()
